# Refactors corrected offers generation
# - Rename header columns for clarity
# - Retain promotional price from consolidated file (update E values)
# - Normalize/rename "Seção" labels for ALTO GIRO rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renaming (row 1) ---
$ws.Range("C1").Value = "ID"
$ws.Range("D1").Value = "PRODUTO"
$ws.Range("E1").Value = "PROMOÇÃO"

# --- Update promotional prices (column E) ---
$priceUpdates = @{
    3  = 16.99
    13 = 11.89
    14 = 14.39
    21 = 17.99
    37 = 2.19
    38 = 2.19
    39 = 2.19
    40 = 2.19
    41 = 2.19
    42 = 2.19
    43 = 2.19
    44 = 2.19
    50 = 5.99
    52 = 12.49
    55 = 2.49
    57 = 8.789999999999999
    60 = 10.19
    61 = 10.19
    64 = 14.99
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $priceUpdates[$row]
}

# --- Rename Seção labels (column B) for ALTO GIRO rows ---
$ws.Range("B29").Value = "#01 MERCEARIA - #02 ALTO GIRO"
$ws.Range("B30").Value = "#01 MERCEARIA - #01 ALTO GIRO"
$ws.Range("B31").Value = "#01 MERCEARIA - #01 ALTO GIRO"
$ws.Range("B32").Value = "#01 MERCEARIA - #01 ALTO GIRO"
$ws.Range("B33").Value = "#01 MERCEARIA - #01 ALTO GIRO"
$ws.Range("B47").Value = "#01 MERCEARIA - #02 ALTO GIRO"
$ws.Range("B56").Value = "#01 MERCEARIA - #02 ALTO GIRO"
